# edit.ps1 - apply the "Finalizing dataset for sharing" changes:
#   1. Bump the cached "datetimeFigureOut" auto-date field shown on the
#      slide master and every slide layout from 7/16/15 -> 7/20/15.
#   2. Center-align the "Resources with Correct RRID" caption paragraph
#      on the single content slide (its sibling captions are already
#      center-aligned).
#
# (The chart's internal axis-id linkage values are an opaque, non-user
# -facing implementation detail that PowerPoint itself regenerates only
# when it rebuilds a chart part from scratch; there is no object-model
# surface to target them directly, so they are intentionally left alone.)

$p = $ppt.ActivePresentation

$oldDate = "7/16/15"
$newDate = "7/20/15"

function Update-DateShape {
    param($shapes)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shape = $shapes.Item($i)
        if ($shape.HasTextFrame) {
            $tr = $shape.TextFrame.TextRange
            if ($tr.Text -eq $oldDate) {
                $tr.Text = $newDate
            }
        }
    }
}

# 1a. Slide master's date placeholder.
$master = $p.SlideMaster
Update-DateShape $master.Shapes

# 1b. Every slide layout's date placeholder.
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DateShape $layout.Shapes
}

# 2. Center the "Resources with Correct RRID" textbox paragraph on slide 1.
$slide = $p.Slides.Item(1)
for ($si = 1; $si -le $slide.Shapes.Count; $si++) {
    $shape = $slide.Shapes.Item($si)
    if ($shape.HasTextFrame) {
        $tr = $shape.TextFrame.TextRange
        if ($tr.Text -like "Resources with Correct RRID*") {
            $para = $tr.Paragraphs(1, 1)
            $para.ParagraphFormat.Alignment = 2  # ppAlignCenter
        }
    }
}
